$wb = $excel.ActiveWorkbook

# --- Sheet1: "Items" (was Sheet1) ---
$items = $wb.Worksheets.Item("Sheet1")
$items.Name = "Items"

# Insert a new blank column at C (this pushes the old product_type/category_id
# columns from C/D to D/E).
$items.Range("C1").EntireColumn.Insert()
$items.Columns.Item(3).ColumnWidth = 12.2

# The old category_id cell (now E3, previously D3) already carries the Consolas
# "code" style (style index 1) - copy its format into the new category column
# before we delete the now-orphaned product_type/category_id columns.
$items.Range("E3").Copy()
$items.Range("C2:C3").PasteSpecial(-4122)
$items.Application.CutCopyMode = $false

# Remove the old product_type (now D) and category_id (now E) columns.
$items.Range("D1:E1").EntireColumn.Delete()

# Header row
$items.Range("A1").Value = "Product No"
$items.Range("B1").Value = "Name"
$items.Range("C1").Value = "category"
$items.Range("D1").Value = "price"
$items.Range("E1").Value = "ingredients"
$items.Range("F1").Value = "country_of_origin"
$items.Range("G1").Value = "year_of_production"
$items.Range("H1").Value = "type_of_drink"
$items.Range("I1").Value = "description"
$items.Range("J1").Value = "is_available" + [char]10 + "(Item is available or not)"
$items.Range("K1").Value = "is_featured" + [char]10 + "(is feature for feature list)"
$items.Range("L1").Value = "is_variable" + [char]10 + "(is variable for variations" + [char]10 + "You can add variations data in Variation Sheet)"

# Row 2 (Product 1)
$items.Range("B2").Value = "New Drink"
$items.Range("C2").Value = "Spirits"
$items.Range("D2").Value = 120

# Row 3 (Product 2)
$items.Range("B3").Value = "New Drink 2"
$items.Range("C3").Value = "Spirits"
$items.Range("D3").Value = 110

# Wrap-text header cells for the boolean columns
$items.Range("J1:L1").WrapText = $true

# Row height + selection for the Items sheet
$items.Rows.Item(1).RowHeight = 165
$items.Range("G17").Select()

# --- Sheet2: "Variations" (was Sheet2) ---
$vars = $wb.Worksheets.Item("Sheet2")
$vars.Name = "Variations"
$vars.Range("G5").Select()
$vars.Activate()

Write-Host "done"
